$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.23512655154714679,
    -0.0059999999757991418,
    -0.0039999999781716866,
    -0.0079999999601199079,
    -0.0029999999768106633,
    -0.0019999999743998131,
    -0.0099999999436461984,
    -0.0099999999437394571,
    -0.0019999999753022024,
    -0.0019999999762774223,
    -0.0029999999722560844,
    -0.0034999999706819551,
    -0.0034999999729281583,
    -0.0079999999565840696,
    -0.00099999998477340313,
    0.028114300870500042,
    -0.0019999999831590287,
    0.0093966756540320162,
    -0.0039999999829967159,
    -0.003999999977695623,
    -0.0039999999764726013,
    -0.043365022910090367,
    -0.0049999999730996336,
    -0.019999999913111743,
    -0.019999999911950894,
    -0.0024999999715085153,
    -0.0024999999696957431,
    -0.0019999999641191479,
    -0.0069999999394685375,
    -0.059999999732885723,
    -0.0069999999349175113,
    -0.0099999999231012993,
    -0.0039999999460693658
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
